$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 289.6930002556438
$ws.Range("G2").Value = 14.51697500019961
$ws.Range("H2").Value = 569.0608923223175
$ws.Range("I2").Value = 0.3442065046112412
$ws.Range("J2").Value = 0.0006744794686918514
$ws.Range("K2").Value = 0.789302625848648
$ws.Range("L2").Value = 0.1882169879485058
$ws.Range("M2").Value = 0.00546515762386496
$ws.Range("N2").Value = 0.3847699617732647

# Row 3
$ws.Range("F3").Value = 0.003186089171715058
$ws.Range("G3").Value = 0.00170688523894819
$ws.Range("H3").Value = 0.004712686832480318
$ws.Range("I3").Value = 0.002943082152674318
$ws.Range("J3").Value = 0.001567559644537148
$ws.Range("K3").Value = 0.004359969912517357
$ws.Range("L3").Value = 0.003300592669246733
$ws.Range("M3").Value = 0.001805281514600336
$ws.Range("N3").Value = 0.004843888341234185

# Row 4
$ws.Range("F4").Value = 289.6961863448155
$ws.Range("G4").Value = 14.51868188543856
$ws.Range("H4").Value = 569.06560500915
$ws.Range("I4").Value = 0.3471495867639156
$ws.Range("J4").Value = 0.002242039113228999
$ws.Range("K4").Value = 0.7936625957611653
$ws.Range("L4").Value = 0.1915175806177525
$ws.Range("M4").Value = 0.007270439138465296
$ws.Range("N4").Value = 0.3896138501144989
